{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the Review_316 diff: updates the header date/title line, rewrites\n// the body paragraphs with the new \"Contextual Document Embeddings\" review\n// text, inserts the additional new paragraphs before the arXiv link, and\n// swaps the arXiv URL for the new paper.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---- 1) Header paragraph: date + title ---------------------------------\n// (paragraph 1 holds two runs split by a manual line break <w:br/>, so we\n// target each run's text individually via search/replace)\n{\n  const dateResults = body.search(\"\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -10.10.24: \u26a1\ufe0f\ud83d\ude80\", { matchCase: true });\n  dateResults.load(\"items\");\n  await context.sync();\n  dateResults.items[0].insertText(\"\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -08.10.24: \u26a1\ufe0f\ud83d\ude80\", \"Replace\");\n  await context.sync();\n\n  const titleResults = body.search(\"DIFFERENTIAL TRANSFORMER\", { matchCase: true });\n  titleResults.load(\"items\");\n  await context.sync();\n  titleResults.items[0].insertText(\"CONTEXTUAL DOCUMENT EMBEDDINGS\", \"Replace\");\n  await context.sync();\n}\n\n// ---- 2) Paragraphs 2-5: replace the whole paragraph text ---------------\n// clear() + insertText(..., \"End\") swaps the paragraph's run content\n// cleanly (avoids carrying over a stale xml:space=\"preserve\" from the\n// original run when the new text has no leading/trailing whitespace).\n{\n  const p2 = paragraphs.items[1];\n  p2.clear();\n  await context.sync();\n  p2.insertText(\"\u05de\u05d6\u05de\u05df \u05dc\u05d0 \u05e1\u05e7\u05e8\u05ea\u05d9 \u05de\u05d0\u05de\u05e8 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05e9\u05dc Document Retrieval \u05d0\u05d5 DG. \u05dc\u05de\u05e2\u05e9\u05d4 DG \u05de\u05d4\u05d5\u05d5\u05d4 \u05e9\u05dc\u05d1 \u05e9\u05dc Retrieval Augmented Generated \u05d0\u05d5 RAG \u05e9\u05de\u05d8\u05e8\u05ea\u05d5 \u05d4\u05d9\u05d0 \u05dc\u05d0\u05ea\u05e8 \u05d0\u05ea \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d4\u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d9\u05dd \u05de\u05e1\u05d8 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd D. \u05d1\u05d3\u05e8\u05da \u05db\u05dc\u05dc \u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05e2\u05dc \u05e1\u05de\u05da \u05e7\u05d9\u05e8\u05d5\u05d1 \u05e9\u05dc \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1(\u05d4\u05e0\u05de\u05d3\u05d3 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05e8\u05d7\u05e7 \u05e7\u05d5\u05e1\u05d9\u05d9\u05df) \u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d5\u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05d4\u05de\u05d5\u05e4\u05e7\u05d9\u05dd \u05e2\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05db\u05dc\u05e9\u05d4\u05d5.\", \"End\");\n  await context.sync();\n}\n\n{\n  const p3 = paragraphs.items[2];\n  p3.clear();\n  await context.sync();\n  p3.insertText(\"\u05d9\u05e9 \u05e9\u05db\u05dc\u05d5\u05dc\u05d9\u05dd \u05dc\u05de\u05e2\u05d8\u05d9\u05dd \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d5, \u05dc\u05de\u05e9\u05dc \u05dc\u05d7\u05dc\u05e7 \u05db\u05dc \u05de\u05e1\u05e4\u05e8 \u05dc\u05e6'\u05d0\u05e0\u05e7\u05d9\u05dd \u05d5\u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc\u05d4\u05dd \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d4\u05e7\u05e8\u05d1\u05d4. \u05d9\u05e6\u05d0 \u05dc\u05d0 \u05de\u05d6\u05de\u05df \u05de\u05d0\u05de\u05e8 \u05e9\u05d4\u05e6\u05d9\u05e2 \u05dc\u05d4\u05d5\u05e1\u05d9\u05e3 \u05ea\u05de\u05e6\u05d5\u05ea \u05dc\u05db\u05dc \u05de\u05e1\u05de\u05da \u05d5\u05db\u05de\u05d5\u05d1\u05df \u05e7\u05d9\u05d9\u05de\u05d5\u05ea \u05e2\u05e9\u05e8\u05d5\u05ea \u05d0\u05d5 \u05de\u05d0\u05d5\u05ea \u05d0\u05d7\u05e8\u05d5\u05ea. \", \"End\");\n  await context.sync();\n}\n\n{\n  const p4 = paragraphs.items[3];\n  p4.clear();\n  await context.sync();\n  p4.insertText(\"\u05d0\u05dd \u05d9\u05e9 \u05d1\u05d9\u05d3\u05d9\u05e0\u05d5 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea D_T \u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05de- (\u05e9\u05d0\u05dc\u05d4, \u05de\u05e1\u05de\u05da \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9) \u05d0\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e2\u05e9\u05d5\u05ea \u05e4\u05d9\u05d9\u05e0\u05d8\u05d9\u05d5\u05df \u05dc\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05db\u05d0\u05dc\u05d5, \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0\u05de\u05d5 \u05e9\u05e0\u05d9 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd: \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d5\u05d4\u05e9\u05e0\u05d9 \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4. \u05d1\u05d3\u05f4\u05db \u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05e2\u05dd \u05dc\u05de\u05d9\u05d3\u05d4 \u05e0\u05d9\u05d2\u05d5\u05d3\u05d9\u05ea \u05e9\u05de\u05d0\u05d5\u05de\u05e0\u05ea \u05dc\u05e7\u05e8\u05d1 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05e9\u05dc \u05db\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05de\u05e1\u05de\u05da \u05d4\u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9 \u05dc\u05d5 \u05d5\u05de\u05e8\u05d7\u05d9\u05e7\u05d4 \u05d0\u05d5\u05ea\u05d5 \u05de\u05db\u05dc \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e9\u05d0\u05e8 \u05de\u05e1\u05de\u05db\u05d9\u05dd.\", \"End\");\n  await context.sync();\n}\n\n{\n  const p5 = paragraphs.items[4];\n  p5.clear();\n  await context.sync();\n  p5.insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05e9\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d4\u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d6\u05d4 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d5\u05e1\u05e4\u05ea \u05e7\u05d5\u05e0\u05d8\u05e7\u05e1\u05d8 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd (=\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1) \u05d4\u05d0\u05dc\u05d5. \u05d0\u05dd \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05de\u05e1\u05de\u05da \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05e9\u05d9\u05d9\u05da \u05d0\u05d5\u05ea\u05d5 \u05dc\u05db\u05de\u05d4 \u05ea\u05d7\u05d5\u05de\u05d9\u05dd (=\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd) \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05e9\u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d9\u05e9\u05ea\u05e0\u05d4 \u05d1\u05d4\u05ea\u05d0\u05dd \u05d1\u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea \u05e6\u05e4\u05d5\u05d9\u05d5\u05ea \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d4\u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05e9\u05dc \u05e8\u05e4\u05d5\u05d0\u05d4 \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05e9\u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05d9\u05e9\u05e7\u05e4\u05d5 \u05d0\u05ea \u05d4\u05d0\u05e1\u05e4\u05e7\u05d8\u05d9\u05dd \u05d4\u05e8\u05e4\u05d5\u05d0\u05d9\u05d9\u05dd \u05d5\u05e2\u05d1\u05d5\u05e8 \u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05d4\u05e1\u05e4\u05d5\u05e8\u05d8 \u05e9\u05d9\u05d4\u05d9\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05f4\u05de\u05db\u05d5\u05d5\u05df\u05f4 \u05dc\u05e1\u05e4\u05d5\u05e8\u05d8. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05db\u05d0\u05df contextualized embedding \u05d1\u05ea\u05dc\u05d5\u05ea \u05d1\u05e9\u05d0\u05dc\u05d5\u05ea \u05de-D_T \u05d5\u05d1\u05e1\u05d8 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd D \u05d1\u05e2\u05e6\u05de\u05d5.\", \"End\");\n  await context.sync();\n}\n\n// ---- 3) Insert the new paragraphs before the link, then update the link ---\n{\n  const linkResults = body.search(\"https://arxiv.org/abs/2410.05258\", { matchCase: true });\n  linkResults.load(\"items\");\n  await context.sync();\n\n  const urlParagraph = linkResults.items[0].paragraphs.getFirst();\n  await context.sync();\n\n  let p = urlParagraph.insertParagraph(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05e8 \u05dc\u05e2\u05e9\u05d5\u05ea \u05d6\u05d0\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc\u05d9 embedding \u05dc\u05de\u05e1\u05de\u05da \u05d0\u05d5 \u05dc\u05d8\u05e7\u05e1\u05d8 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d4\u05d1\u05d0\u05d4. \u05e7\u05d5\u05d3\u05dd \u05db\u05dc \u05d0\u05e0\u05d5 \u05de\u05d7\u05dc\u05e7\u05d9\u05dd \u05d0\u05ea D \u05dc\u05db\u05de\u05d4 \u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05dc\u05e4\u05d9 \u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd (\u05e2\u05dd \u05de\u05d5\u05d3\u05dc embedding \u05d4\u05ea\u05d7\u05dc\u05ea\u05d9). \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05de\u05e7\u05e1\u05de\u05d9\u05dd \u05d0\u05ea \u05e1\u05db\u05d5\u05de\u05d9 \u05d4\u05dc\u05d5\u05e1\u05d9\u05dd \u05d4\u05e0\u05d9\u05d2\u05d5\u05d3\u05d9\u05d9\u05dd \u05e2\u05dc \u05e4\u05e0\u05d9 \u05db\u05dc \u05d4\u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d1\u05e0\u05d5\u05ea \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d4 \u05d5\u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05da \u05db\u05da \u05e9:\", \"Before\");\n  await context.sync();\n  p = p.insertParagraph(\"\u05f4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05d5\u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05da \u05d4\u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9 \u05dc\u05d4 \u05d9\u05d4\u05d9\u05d5 \u05e7\u05e8\u05d5\u05d1\u05d9\u05dd \u05d0\u05d7\u05d3 \u05dc\u05e9\u05e0\u05d9 \u05d1\u05ea\u05d5\u05da \u05db\u05dc \u05e7\u05dc\u05e1\u05d8\u05e8 (\u05d4\u05de\u05d3\u05de\u05d4 \u05d3\u05d5\u05de\u05d9\u05d9\u05df) \u05d5\u05d0\u05d9\u05dc\u05d5 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05d9\u05d4\u05d9\u05d4 \u05e8\u05d7\u05d5\u05e7 \u05de\u05d4\u05db\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d9\u05dd \u05d1\u05e7\u05dc\u05e1\u05d8\u05e8\u05f4. \", \"After\");\n  await context.sync();\n  p = p.insertParagraph(\"\u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05de\u05ea\u05d0\u05d9\u05de\u05d9\u05dd \u05d0\u05ea \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05db\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea. \u05d4\u05de\u05d0\u05de\u05e8 \u05d2\u05dd \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d1\u05e0\u05d9\u05d4 \u05e9\u05dc \u05d1\u05d0\u05e6'\u05d9\u05dd (\u05db\u05db\u05d4 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e8\u05e9\u05ea\u05d5\u05ea \u05d4\u05d9\u05d5\u05dd) \u05db\u05da \u05e9\u05d4\u05e8\u05e9\u05ea \u05ea\u05dc\u05de\u05d3 \u05e2\u05dc \u05e9\u05d9\u05dc\u05d5\u05d1\u05d9 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d4\u05e7\u05e9\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8(\u05dc\u05de\u05e9\u05dc \u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d3\u05d5\u05de\u05d9\u05dd \u05e1\u05de\u05e0\u05d8\u05d9\u05ea \u05d0\u05d1\u05dc \u05de\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd).\", \"After\");\n  await context.sync();\n  p = p.insertParagraph(\"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05e9\u05dc\u05d1 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05dc\u05d1\u05e0\u05d9\u05d9\u05d4 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05de\u05e1\u05de\u05da \u05e0\u05ea\u05d5\u05df 'D. \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05de\u05e1\u05de\u05da 'D \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05e9\u05e8\u05e9\u05d5\u05e8 \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05db\u05dc  \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d5\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05e9\u05dc \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05de 'D (\u05e9\u05d4\u05dd \u05ea\u05dc\u05d5\u05d9\u05d9 \u05d4\u05e7\u05e9\u05e8 \u05d4\u05de\u05e1\u05de\u05da \u05db\u05de\u05d5\u05d1\u05df). \u05d1\u05d4\u05de\u05e9\u05da \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05dc\u05de\u05e1\u05de\u05da \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d5\u05de\u05d4 \u05dc\u05de\u05d4 \u05e9\u05ea\u05d5\u05d0\u05e8 \u05dc\u05e4\u05e0\u05d9 \u05d0\u05d1\u05dc \u05e2\u05dd \u05db\u05de\u05d4 \u05d8\u05e8\u05d9\u05e7\u05d9\u05dd \u05dc\u05d9\u05d9\u05e2\u05d5\u05dc \u05d4\u05d0\u05d9\u05de\u05d5\u05df.\", \"After\");\n  await context.sync();\n  p = p.insertParagraph(\"\u05d0\u05e6\u05d9\u05d9\u05df \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05db\u05ea\u05d5\u05d1 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d1\u05e8\u05d5\u05e8\u05d4\u2026.\", \"After\");\n  await context.sync();\n\n  linkResults.items[0].insertText(\"https://arxiv.org/abs/2410.02525\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the Review_316 diff: updates the header date/title line, rewrites\n# the body paragraphs with the new \"Contextual Document Embeddings\" review\n# text, inserts the additional new paragraphs before the arXiv link, and\n# swaps the arXiv URL for the new paper.\n\n$d = $word.ActiveDocument\n\nfunction Replace-WholeText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# ---- 1) Header paragraph: date + title ----------------------------------\nReplace-WholeText \"\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -10.10.24: \u26a1\ufe0f\ud83d\ude80\" \"\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -08.10.24: \u26a1\ufe0f\ud83d\ude80\"\nReplace-WholeText \"DIFFERENTIAL TRANSFORMER\" \"CONTEXTUAL DOCUMENT EMBEDDINGS\"\n\n# ---- 2) Paragraphs 2-5: replace the whole paragraph text ----------------\nReplace-WholeText \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e2\u05e9\u05d4 \u05d4\u05e8\u05d1\u05d4 \u05d2\u05dc\u05d9\u05dd \u05d1\u05d9\u05d5\u05de\u05d9\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d9\u05dd \u05d5\u05d6\u05d5 \u05d4\u05e1\u05d9\u05d1\u05d4 \u05e9\u05d1\u05d7\u05e8\u05ea\u05d9 \u05d0\u05d5\u05ea\u05d5 \u05dc\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d4\u05d9\u05d5\u05de\u05d9\u05ea \u05e9\u05dc\u05d9. \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d7\u05d6\u05d9\u05e8 \u05d0\u05d5\u05ea\u05d9 3-4 \u05e9\u05e0\u05d9\u05dd \u05d0\u05d7\u05d5\u05e8\u05d4 \u05dc\u05ea\u05e7\u05d5\u05e4\u05d4 \u05e9\u05d1\u05d4 \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 \u05d9\u05de\u05d9 \u05d9\u05e6\u05d0\u05d5 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05d4\u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05db\u05dc\u05d5\u05dc\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05dc\u05dc\u05d9\u05d1\u05d4 \u05e9\u05dc \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05db\u05d4 \u05d0\u05d4\u05d5\u05d1\u05d9\u05dd \u05e2\u05dc\u05d9\u05e0\u05d5. \u05db\u05de\u05d5\u05d1\u05df \u05d0\u05e0\u05d9 \u05de\u05ea\u05db\u05d5\u05d5\u05df \u05dc\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d4-attention \u05e9\u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e0\u05d5 \u05dc\u05db\u05de\u05ea \u05e7\u05e9\u05e8\u05d9\u05dd \u05d1\u05d9\u05df \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05d1\u05d8\u05e7\u05e1\u05d8. \" \"\u05de\u05d6\u05de\u05df \u05dc\u05d0 \u05e1\u05e7\u05e8\u05ea\u05d9 \u05de\u05d0\u05de\u05e8 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05e9\u05dc Document Retrieval \u05d0\u05d5 DG. \u05dc\u05de\u05e2\u05e9\u05d4 DG \u05de\u05d4\u05d5\u05d5\u05d4 \u05e9\u05dc\u05d1 \u05e9\u05dc Retrieval Augmented Generated \u05d0\u05d5 RAG \u05e9\u05de\u05d8\u05e8\u05ea\u05d5 \u05d4\u05d9\u05d0 \u05dc\u05d0\u05ea\u05e8 \u05d0\u05ea \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d4\u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d9\u05dd \u05de\u05e1\u05d8 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd D. \u05d1\u05d3\u05e8\u05da \u05db\u05dc\u05dc \u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05e2\u05dc \u05e1\u05de\u05da \u05e7\u05d9\u05e8\u05d5\u05d1 \u05e9\u05dc \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1(\u05d4\u05e0\u05de\u05d3\u05d3 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05e8\u05d7\u05e7 \u05e7\u05d5\u05e1\u05d9\u05d9\u05df) \u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d5\u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05d4\u05de\u05d5\u05e4\u05e7\u05d9\u05dd \u05e2\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05db\u05dc\u05e9\u05d4\u05d5.\"\nReplace-WholeText \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05e6\u05d9\u05e2\u05d5 \u05dc\u05d4\u05d7\u05dc\u05d9\u05e3 \u05d0\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1 \u05d4\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d4\u05e8\u05d2\u05d9\u05dc \u05e9\u05d9\u05e9 \u05dc\u05e0\u05d5 \u05d1\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd \u05d1\u05d4\u05e4\u05e8\u05e9 \u05de\u05e9\u05d5\u05e7\u05dc\u05dc (\u05e8\u05e7 \u05d4\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d4\u05e9\u05e0\u05d9 \u05de\u05e9\u05d5\u05e7\u05dc\u05dc) \u05e9\u05dc \u05d4\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1\u05d9\u05dd. \u05db\u05dc \u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05de\u05d7\u05d5\u05e9\u05d1 \u05e2\u05dd \u05de\u05d8\u05e8\u05d9\u05e6\u05ea Q \u05d5-K \u05de\u05e9\u05dc\u05d4 \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05e9\u05e7\u05d5\u05dc \u03bb \u05e9\u05dc \u05d4\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d4\u05e9\u05e0\u05d9 \u05de\u05d7\u05d5\u05e9\u05d1 \u05d1\u05d0\u05d5\u05e4\u05df \u05d4\u05d1\u05d0:  \u03bb = exp(\u03bb_q1 \u00b7 \u03bb_k1 ) \u2212 exp(\u03bb_q2 \u00b7 \u03bb_k2 ) + \u03bb_init \u05db\u05d0\u05e9\u05e8 \" \"\u05d9\u05e9 \u05e9\u05db\u05dc\u05d5\u05dc\u05d9\u05dd \u05dc\u05de\u05e2\u05d8\u05d9\u05dd \u05dc\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d6\u05d5, \u05dc\u05de\u05e9\u05dc \u05dc\u05d7\u05dc\u05e7 \u05db\u05dc \u05de\u05e1\u05e4\u05e8 \u05dc\u05e6'\u05d0\u05e0\u05e7\u05d9\u05dd \u05d5\u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1\u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc\u05d4\u05dd \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d4\u05e7\u05e8\u05d1\u05d4. \u05d9\u05e6\u05d0 \u05dc\u05d0 \u05de\u05d6\u05de\u05df \u05de\u05d0\u05de\u05e8 \u05e9\u05d4\u05e6\u05d9\u05e2 \u05dc\u05d4\u05d5\u05e1\u05d9\u05e3 \u05ea\u05de\u05e6\u05d5\u05ea \u05dc\u05db\u05dc \u05de\u05e1\u05de\u05da \u05d5\u05db\u05de\u05d5\u05d1\u05df \u05e7\u05d9\u05d9\u05de\u05d5\u05ea \u05e2\u05e9\u05e8\u05d5\u05ea \u05d0\u05d5 \u05de\u05d0\u05d5\u05ea \u05d0\u05d7\u05e8\u05d5\u05ea. \"\nReplace-WholeText \"\u03bb_q1 , \u03bb_k1 , \u03bb_q2 , \u03bb_k2 \u2208  R^d \u05d4\u05d9\u05e0\u05dd \u05e0\u05dc\u05de\u05d3\u05d9\u05dd \u05d5- ((\u03bb_init = 0.8 \u2212 0.6 \u00d7 exp(\u22120.3 \u00b7 (l \u2212 1, \u05db\u05d0\u05e9\u05e8 l \u05d6\u05d4 \u05de\u05e1\u05e4\u05e8 \u05d4\u05e9\u05db\u05d1\u05d4 (\u05e9\u05dc \u05d1\u05dc\u05d5\u05e7 \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8). \u05d0\u05dd \u05d4\u05e0\u05d5\u05e1\u05d7\u05d4 \u05e2\u05d1\u05d5\u05e8 \u03bb \u05d0\u05d9\u05db\u05e9\u05d4\u05d5 \u05de\u05d5\u05d1\u05e0\u05ea \u05d5\u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea \u05d4\u05e0\u05d5\u05e1\u05d7\u05d4 \u05e2\u05d1\u05d5\u05e8 \u03bb_init \u05e0\u05d5\u05ea\u05e8\u05ea \u05d1\u05d2\u05d3\u05e8 \u05ea\u05e2\u05dc\u05d5\u05de\u05d4 (\u05d0\u05dc\u05d0 \u05d0\u05dd \u05db\u05df \u05d6\u05d4 \u05e0\u05d9\u05e1\u05d5\u05d9 \u05d5\u05ea\u05d4\u05d9\u05d4 \u05e8\u05d2\u05e8\u05e1\u05d9\u05d4 \u05e9\u05dc \u05d4\u05e2\u05e8\u05db\u05d9\u05dd \u05e9\u05d4\u05ea\u05e7\u05d1\u05dc\u05d5 \u05e2\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05de\u05e6\u05d5\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea).\" \"\u05d0\u05dd \u05d9\u05e9 \u05d1\u05d9\u05d3\u05d9\u05e0\u05d5 \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea D_T \u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d9\u05dd \u05de- (\u05e9\u05d0\u05dc\u05d4, \u05de\u05e1\u05de\u05da \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9) \u05d0\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e2\u05e9\u05d5\u05ea \u05e4\u05d9\u05d9\u05e0\u05d8\u05d9\u05d5\u05df \u05dc\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05db\u05d0\u05dc\u05d5, \u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0\u05de\u05d5 \u05e9\u05e0\u05d9 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd: \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d5\u05d4\u05e9\u05e0\u05d9 \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4. \u05d1\u05d3\u05f4\u05db \u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05e2\u05dd \u05dc\u05de\u05d9\u05d3\u05d4 \u05e0\u05d9\u05d2\u05d5\u05d3\u05d9\u05ea \u05e9\u05de\u05d0\u05d5\u05de\u05e0\u05ea \u05dc\u05e7\u05e8\u05d1 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05e9\u05dc \u05db\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d4\u05de\u05e1\u05de\u05da \u05d4\u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9 \u05dc\u05d5 \u05d5\u05de\u05e8\u05d7\u05d9\u05e7\u05d4 \u05d0\u05d5\u05ea\u05d5 \u05de\u05db\u05dc \u05de\u05d4\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd \u05e9\u05dc \u05e9\u05d0\u05e8 \u05de\u05e1\u05de\u05db\u05d9\u05dd.\"\nReplace-WholeText \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05ea\u05d5\u05e6\u05d0\u05ea \u05d0\u05d1\u05dc \u05d4\u05d1\u05d3\u05d9\u05e7\u05d5\u05ea \u05e0\u05e2\u05e9\u05d5 \u05d1\u05e2\u05d9\u05e7\u05e8 \u05dc\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd 3B \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd. \u05d9\u05e9 \u05d2\u05dd \u05d8\u05e2\u05e0\u05d5\u05ea \u05dc\u05e7\u05e0\u05e1\u05d5\u05dc \u05e9\u05dc \u05e8\u05e2\u05e9 \u05db\u05dc\u05e9\u05d4\u05d5 \u05e9\u05d0\u05e0\u05d9 \u05dc\u05d0 \u05d1\u05d8\u05d5\u05d7 \u05e9\u05d0\u05e0\u05d9 \u05de\u05d1\u05d9\u05df. \u05d1\u05e7\u05d9\u05e6\u05e8 \u05d0\u05e0\u05d9 \u05e7\u05e6\u05ea \u05e1\u05e7\u05e4\u05d8\u05d9, \u05de\u05d5\u05d3\u05d4\u2026.\" \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05e9\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d4\u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d6\u05d4 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d5\u05e1\u05e4\u05ea \u05e7\u05d5\u05e0\u05d8\u05e7\u05e1\u05d8 \u05dc\u05d9\u05d9\u05e6\u05d5\u05d2\u05d9\u05dd (=\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1) \u05d4\u05d0\u05dc\u05d5. \u05d0\u05dd \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05de\u05e1\u05de\u05da \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05e9\u05d9\u05d9\u05da \u05d0\u05d5\u05ea\u05d5 \u05dc\u05db\u05de\u05d4 \u05ea\u05d7\u05d5\u05de\u05d9\u05dd (=\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd) \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05e9\u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d9\u05e9\u05ea\u05e0\u05d4 \u05d1\u05d4\u05ea\u05d0\u05dd \u05d1\u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05dd \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea \u05e6\u05e4\u05d5\u05d9\u05d5\u05ea \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05d4\u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05e9\u05dc \u05e8\u05e4\u05d5\u05d0\u05d4 \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05e9\u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05d9\u05e9\u05e7\u05e4\u05d5 \u05d0\u05ea \u05d4\u05d0\u05e1\u05e4\u05e7\u05d8\u05d9\u05dd \u05d4\u05e8\u05e4\u05d5\u05d0\u05d9\u05d9\u05dd \u05d5\u05e2\u05d1\u05d5\u05e8 \u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05d4\u05e1\u05e4\u05d5\u05e8\u05d8 \u05e9\u05d9\u05d4\u05d9\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05f4\u05de\u05db\u05d5\u05d5\u05df\u05f4 \u05dc\u05e1\u05e4\u05d5\u05e8\u05d8. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05db\u05d0\u05df contextualized embedding \u05d1\u05ea\u05dc\u05d5\u05ea \u05d1\u05e9\u05d0\u05dc\u05d5\u05ea \u05de-D_T \u05d5\u05d1\u05e1\u05d8 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd D \u05d1\u05e2\u05e6\u05de\u05d5.\"\n\n# ---- 3) Insert the new paragraphs before the link, then update the link --\n$urlParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$urlParagraph.Range.InsertBefore(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05e8 \u05dc\u05e2\u05e9\u05d5\u05ea \u05d6\u05d0\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d0\u05d9\u05de\u05d5\u05df \u05de\u05d5\u05d3\u05dc\u05d9 embedding \u05dc\u05de\u05e1\u05de\u05da \u05d0\u05d5 \u05dc\u05d8\u05e7\u05e1\u05d8 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d4\u05d1\u05d0\u05d4. \u05e7\u05d5\u05d3\u05dd \u05db\u05dc \u05d0\u05e0\u05d5 \u05de\u05d7\u05dc\u05e7\u05d9\u05dd \u05d0\u05ea D \u05dc\u05db\u05de\u05d4 \u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05dc\u05e4\u05d9 \u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd (\u05e2\u05dd \u05de\u05d5\u05d3\u05dc embedding \u05d4\u05ea\u05d7\u05dc\u05ea\u05d9). \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05de\u05e7\u05e1\u05de\u05d9\u05dd \u05d0\u05ea \u05e1\u05db\u05d5\u05de\u05d9 \u05d4\u05dc\u05d5\u05e1\u05d9\u05dd \u05d4\u05e0\u05d9\u05d2\u05d5\u05d3\u05d9\u05d9\u05dd \u05e2\u05dc \u05e4\u05e0\u05d9 \u05db\u05dc \u05d4\u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05d4\u05d0\u05dc\u05d5. \u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05e8\u05d5\u05e6\u05d9\u05dd \u05dc\u05d1\u05e0\u05d5\u05ea \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d4 \u05d5\u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05da \u05db\u05da \u05e9:\" + \"`r\")\n\n$urlParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$urlParagraph.Range.InsertBefore(\"\u05f4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05d5\u05e9\u05dc \u05d4\u05de\u05e1\u05de\u05da \u05d4\u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9 \u05dc\u05d4 \u05d9\u05d4\u05d9\u05d5 \u05e7\u05e8\u05d5\u05d1\u05d9\u05dd \u05d0\u05d7\u05d3 \u05dc\u05e9\u05e0\u05d9 \u05d1\u05ea\u05d5\u05da \u05db\u05dc \u05e7\u05dc\u05e1\u05d8\u05e8 (\u05d4\u05de\u05d3\u05de\u05d4 \u05d3\u05d5\u05de\u05d9\u05d9\u05df) \u05d5\u05d0\u05d9\u05dc\u05d5 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05d4\u05e9\u05d0\u05dc\u05d4 \u05d9\u05d4\u05d9\u05d4 \u05e8\u05d7\u05d5\u05e7 \u05de\u05d4\u05db\u05dc \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d9\u05dd \u05d1\u05e7\u05dc\u05e1\u05d8\u05e8\u05f4. \" + \"`r\")\n\n$urlParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$urlParagraph.Range.InsertBefore(\"\u05db\u05dc\u05d5\u05de\u05e8 \u05d0\u05e0\u05d5 \u05de\u05ea\u05d0\u05d9\u05de\u05d9\u05dd \u05d0\u05ea \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05db\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d3\u05d5\u05de\u05d9\u05d9\u05df \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea. \u05d4\u05de\u05d0\u05de\u05e8 \u05d2\u05dd \u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d1\u05e0\u05d9\u05d4 \u05e9\u05dc \u05d1\u05d0\u05e6'\u05d9\u05dd (\u05db\u05db\u05d4 \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05e8\u05e9\u05ea\u05d5\u05ea \u05d4\u05d9\u05d5\u05dd) \u05db\u05da \u05e9\u05d4\u05e8\u05e9\u05ea \u05ea\u05dc\u05de\u05d3 \u05e2\u05dc \u05e9\u05d9\u05dc\u05d5\u05d1\u05d9 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d4\u05e7\u05e9\u05d9\u05dd \u05d1\u05d9\u05d5\u05ea\u05e8(\u05dc\u05de\u05e9\u05dc \u05de\u05e1\u05de\u05db\u05d9\u05dd \u05d3\u05d5\u05de\u05d9\u05dd \u05e1\u05de\u05e0\u05d8\u05d9\u05ea \u05d0\u05d1\u05dc \u05de\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd).\" + \"`r\")\n\n$urlParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$urlParagraph.Range.InsertBefore(\"\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05e9\u05dc\u05d1 \u05d0\u05ea \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05dc\u05d1\u05e0\u05d9\u05d9\u05d4 \u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc \u05de\u05e1\u05de\u05da \u05e0\u05ea\u05d5\u05df 'D. \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05d9\u05e6\u05d5\u05d2 \u05e9\u05dc \u05de\u05e1\u05de\u05da 'D \u05de\u05d5\u05e8\u05db\u05d1 \u05de\u05e9\u05e8\u05e9\u05d5\u05e8 \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05db\u05dc  \u05d4\u05de\u05e1\u05de\u05db\u05d9\u05dd \u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d5\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1 \u05e9\u05dc \u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05de 'D (\u05e9\u05d4\u05dd \u05ea\u05dc\u05d5\u05d9\u05d9 \u05d4\u05e7\u05e9\u05e8 \u05d4\u05de\u05e1\u05de\u05da \u05db\u05de\u05d5\u05d1\u05df). \u05d1\u05d4\u05de\u05e9\u05da \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05e0\u05e7\u05d5\u05d3\u05e8 \u05dc\u05de\u05e1\u05de\u05da \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d5\u05de\u05d4 \u05dc\u05de\u05d4 \u05e9\u05ea\u05d5\u05d0\u05e8 \u05dc\u05e4\u05e0\u05d9 \u05d0\u05d1\u05dc \u05e2\u05dd \u05db\u05de\u05d4 \u05d8\u05e8\u05d9\u05e7\u05d9\u05dd \u05dc\u05d9\u05d9\u05e2\u05d5\u05dc \u05d4\u05d0\u05d9\u05de\u05d5\u05df.\" + \"`r\")\n\n$urlParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$urlParagraph.Range.InsertBefore(\"\u05d0\u05e6\u05d9\u05d9\u05df \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05dc\u05d0 \u05db\u05ea\u05d5\u05d1 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d1\u05e8\u05d5\u05e8\u05d4\u2026.\" + \"`r\")\n\nReplace-WholeText \"https://arxiv.org/abs/2410.05258\" \"https://arxiv.org/abs/2410.02525\"\n"}
